$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 36 (pushes existing rows 36..140 down to 37..141)
$ws.Rows("36:36").Insert()

# Populate the newly inserted row 36 with the new data record
$ws.Cells.Item(36, 1).Value = 10
$ws.Cells.Item(36, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(36, 3).Value = "La Araucanía"
$ws.Cells.Item(36, 4).Value = 44998
$ws.Cells.Item(36, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 5).Value = 9
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100101
$ws.Cells.Item(36, 8).Value = "Berries"
$ws.Cells.Item(36, 9).Value = 100101001
$ws.Cells.Item(36, 10).Value = "Arándano (blue)"
$ws.Cells.Item(36, 11).Value = "Sin especificar"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 200
$ws.Cells.Item(36, 14).Value = 1500
$ws.Cells.Item(36, 15).Value = 1500
$ws.Cells.Item(36, 16).Value = 1500
$ws.Cells.Item(36, 17).Value = "$/kilo"
$ws.Cells.Item(36, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(36, 19).Value = 1500
$ws.Cells.Item(36, 20).Value = 1
